# Habitat wetlands workbook update:
#  - sort the Sheet1 data table (rows 2:50) by rgn_id (A) ascending, then
#    year (C) ascending - this is what breaks up the old shared formulas
#    in column E and leaves the sheet's sortState behind.
#  - the row that ends up at row 16 (Maui Nui / Estuarine Emergent Wetland,
#    2005) gets its computed ratio hard-coded to 0.89 and a note added in
#    column I referencing the source of that override.
#  - tidy up the selections that were left on Sheet1 / Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Use the worksheet Sort object (rather than Range.Sort) so the applied
# sort is remembered on the sheet (sortState), same as Excel's Data > Sort
# dialog leaves behind.
$sortObj = $ws1.Sort
$sortFields = $sortObj.SortFields
$sortFields.Clear()
$sortFields.Add($ws1.Range("A2:A50"))
$sortFields.Add($ws1.Range("C2:C50"))
$sortObj.SetRange($ws1.Range("A2:I50"))
$sortObj.Header = 2
$sortObj.Apply()

# Override the extent/condition ratio for the row identified in the commit
# message ("data from excel change 2005/2010 tables") with the value taken
# from the 2005/2010 comparison tables, and note where it came from.
$ws1.Range("E16").Value = 0.89
$ws1.Range("I16").Value = "data from excel change 2005/2010 tables"

# Restore/settle the UI selection state: Sheet1 stays the active tab with
# its selection moved to E7, while Sheet2 keeps a B1 selection in reserve.
$ws2.Activate()
$ws2.Range("B1").Select()
$ws1.Activate()
$ws1.Range("E7").Select()
